$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"="1.02"; "C"="1.016575399075721"; "D"="1.022036165938839"; "E"="1.018079921344615"; "F"="1.01492247086263"; "I"="1.026428262455471"; "J"="1.021794426393935"; "K"="1.024871731696536"; "L"="1.020927201768244"; "M"="1.017779172051364" }
    3 = @{ "B"="1.02"; "C"="1.01758499701192"; "D"="1.022738956851938"; "E"="1.01893769224147"; "F"="1.016575723826187"; "I"="1.02655242823776"; "J"="1.022439219096827"; "K"="1.025381386959106"; "L"="1.021590562361069"; "M"="1.019235123611414" }
    4 = @{ "B"="1.019999999999999"; "C"="1.018238117224918"; "D"="1.023193354322733"; "E"="1.019492974428462"; "F"="1.017645160628615"; "I"="1.026631240251872"; "J"="1.022855740859572"; "K"="1.025710143057602"; "L"="1.022019413185946"; "M"="1.020176446575417" }
    5 = @{ "B"="1.02"; "C"="1.018512652216644"; "D"="1.023384297700813"; "E"="1.01972647432373"; "F"="1.018094680237479"; "I"="1.026664006033661"; "J"="1.023030678983413"; "K"="1.025848106794078"; "L"="1.022199609524285"; "M"="1.020571999226914" }
    6 = @{ "B"="1.02"; "C"="1.018558745720429"; "D"="1.023416352922965"; "E"="1.019765683462125"; "F"="1.018170152560844"; "I"="1.026669486040007"; "J"="1.023060042048678"; "K"="1.025871257127548"; "L"="1.022229859863914"; "M"="1.020638403905121" }
    7 = @{ "B"="1.02"; "C"="1.018241785718215"; "D"="1.023195906053109"; "E"="1.019496094232499"; "F"="1.017651167404378"; "I"="1.026631679511256"; "J"="1.022858079049406"; "K"="1.025711987499738"; "L"="1.022021821341536"; "M"="1.020181732667243" }
    8 = @{ "B"="1.02"; "C"="1.016916630497611"; "D"="1.02227375081335"; "E"="1.018369757953472"; "F"="1.015481267661749"; "I"="1.026470541901546"; "J"="1.022012482806205"; "K"="1.025044184327718"; "L"="1.021151468073679"; "M"="1.018271381809822" }
    9 = @{ "B"="1.02"; "C"="1.014580293964472"; "D"="1.020646085865126"; "E"="1.016386895540633"; "F"="1.011654798108463"; "I"="1.026174872547915"; "J"="1.02051703326557"; "K"="1.023859577264158"; "L"="1.019614809869173"; "M"="1.014898893402669" }
    10 = @{ "B"="1.02"; "C"="1.013021841066643"; "D"="1.019559166275683"; "E"="1.015066243486943"; "F"="1.009101476840297"; "I"="1.025969883242132"; "J"="1.019516401444701"; "K"="1.023064553125041"; "L"="1.018588341998398"; "M"="1.012646016316115" }
    11 = @{ "B"="1.02"; "C"="1.012346787360823"; "D"="1.019088090513704"; "E"="1.014494680898912"; "F"="1.007995207076624"; "I"="1.02587925339289"; "J"="1.01908223903227"; "K"="1.022719042840219"; "L"="1.018143383010055"; "M"="1.011669331470965" }
    12 = @{ "B"="1.02"; "C"="1.012096005985239"; "D"="1.018913047064612"; "E"="1.014282420028335"; "F"="1.007584180848626"; "I"="1.025845308785465"; "J"="1.01892083842139"; "K"="1.022590515406041"; "L"="1.017978030978306"; "M"="1.011306363467904" }
    13 = @{ "B"="1.02"; "C"="1.012149801108823"; "D"="1.018950597416832"; "E"="1.014327948770025"; "F"="1.007672352411251"; "I"="1.025852602718548"; "J"="1.018955465448314"; "K"="1.022618093558412"; "L"="1.018013502929098"; "M"="1.01138422977092" }
    14 = @{ "B"="1.02"; "C"="1.012326058443739"; "D"="1.019073622692163"; "E"="1.014477134475953"; "F"="1.007961233780678"; "I"="1.02587645324605"; "J"="1.019068900332674"; "K"="1.022708422590504"; "L"="1.018129716481771"; "M"="1.011639332214027" }
    15 = @{ "B"="1.02"; "C"="1.012434651545663"; "D"="1.019149414020556"; "E"="1.014569058388275"; "F"="1.008139208510017"; "I"="1.02589111115538"; "J"="1.019138773612734"; "K"="1.022764052164874"; "L"="1.018201309603701"; "M"="1.011796484677385" }
    16 = @{ "B"="1.02"; "C"="1.013066637129772"; "D"="1.019590420914886"; "E"="1.015104182255978"; "F"="1.009174881589792"; "I"="1.025975858689142"; "J"="1.019545196717251"; "K"="1.023087456941997"; "L"="1.018617862052851"; "M"="1.012710810306764" }
    17 = @{ "B"="1.02"; "C"="1.013463001596714"; "D"="1.019866937144499"; "E"="1.01543992826175"; "F"="1.009824347977968"; "I"="1.026028518441906"; "J"="1.019799898423503"; "K"="1.023289982784614"; "L"="1.018879022363394"; "M"="1.013284023073536" }
    18 = @{ "B"="1.02"; "C"="1.013694171861702"; "D"="1.020028182660723"; "E"="1.015635791087851"; "F"="1.010203107162639"; "I"="1.026059053736506"; "J"="1.019948376552272"; "K"="1.023407991170436"; "L"="1.01903130535936"; "M"="1.013618255584452" }
    19 = @{ "B"="1.02"; "C"="1.013772991184794"; "D"="1.020083156136956"; "E"="1.015702579988222"; "F"="1.010332243731739"; "I"="1.026069434909006"; "J"="1.019998989352903"; "K"="1.023448208378639"; "L"="1.019083221904484"; "M"="1.013732201324003" }
    20 = @{ "B"="1.02"; "C"="1.013420477769123"; "D"="1.019837273885157"; "E"="1.015403903013792"; "F"="1.009754673016978"; "I"="1.02602288719253"; "J"="1.019772580145495"; "K"="1.023268266248874"; "L"="1.018851007240016"; "M"="1.013222534488014" }
    21 = @{ "B"="1.02"; "C"="1.012274156066409"; "D"="1.019037396607302"; "E"="1.014433201822409"; "F"="1.007876168477913"; "I"="1.025869437607231"; "J"="1.019035500274498"; "K"="1.022681828176307"; "L"="1.018095496547077"; "M"="1.011564216034856" }
    22 = @{ "B"="1.02"; "C"="1.011553207013932"; "D"="1.018534106318378"; "E"="1.013823131584085"; "F"="1.006694446813995"; "I"="1.025771334037443"; "J"="1.018571296915223"; "K"="1.022312014567588"; "L"="1.017620045934155"; "M"="1.01052049852883" }
    23 = @{ "B"="1.02"; "C"="1.011935416200788"; "D"="1.018800945628863"; "E"="1.014146517982557"; "F"="1.007320962519004"; "I"="1.025823494501011"; "J"="1.01881745334963"; "K"="1.02250816381488"; "L"="1.01787213235877"; "M"="1.011073896521112" }
    24 = @{ "B"="1.02"; "C"="1.013439692525261"; "D"="1.019850677563984"; "E"="1.015420181186343"; "F"="1.009786156327183"; "I"="1.026025432268605"; "J"="1.019784924361976"; "K"="1.023278079391863"; "L"="1.018863666214571"; "M"="1.013250318879863" }
    25 = @{ "B"="1.02"; "C"="1.015184446388862"; "D"="1.021067196074201"; "E"="1.016899290678068"; "F"="1.012644412586164"; "I"="1.026252698885156"; "J"="1.020904286338894"; "K"="1.024166757786121"; "L"="1.020012428806697"; "M"="1.015771536848917" }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = [double]$data[$row][$col]
    }
}
